# Remove the trailing "Ver no Jupiter ..." / copyright footer block that
# Jekyll appended after the bibliography entry ending in "...1984.".
#
# Target: delete the blank paragraph, the "Ver no Jupiter Salvar em pdf
# Salvar em docx" paragraph, and the "© 2020 . Contact: ..." paragraph,
# while leaving the remaining blank paragraph and the page-break
# paragraph that follow them untouched.

$d = $word.ActiveDocument

$startText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$endText   = "Contact: luizeleno@usp.br"

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($startPara -eq $null -and $text -like "*$startText*") {
        $startPara = $para
    }
    if ($text -like "*$endText*") {
        $endPara = $para
    }
}

# Include the blank paragraph immediately preceding the "Ver no Jupiter"
# paragraph so it is removed along with the footer text.
$deleteStart = $startPara.Previous().Range.Start
$deleteEnd = $endPara.Range.End

$range = $d.Range($deleteStart, $deleteEnd)
$range.Delete()
